# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# The previous Estado de Cuenta (EC) periods are replaced with the new set
# of periods (most-recent-first), and the "Salario Basico" amount is
# updated for every row of the database.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New "Periodo Mora" (column E) values - previous periods removed, new
# period (1901) added, most recent period first.
$periodos = @("1901", "1812", "1811", "1810", "1809", "1808")

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodos[$i]
}

# Updated "Salario Basico" (column G) for every worker/period row.
$ws.Range("G16:G21").Value = 1423500
